$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 251; shifts existing rows 251:288 down to 252:289
# -4121 = xlShiftDown
$ws.Rows(251).Insert(-4121)

# Populate the newly inserted row 251 with the new record
$ws.Range("A251").Value = 10
$ws.Range("B251").Value = "Vega Modelo de Temuco"
$ws.Range("C251").Value = "La Araucanía"
$ws.Range("D251").Value = 44504
$ws.Range("E251").Value = 9
$ws.Range("F251").Value = "Fruta"
$ws.Range("G251").Value = 100103
$ws.Range("H251").Value = "Frutos de hueso (carozo)"
$ws.Range("I251").Value = 100103006
$ws.Range("J251").Value = "Nectarín"
$ws.Range("K251").Value = "Early Glo"
$ws.Range("L251").Value = "Primera"
$ws.Range("M251").Value = 140
$ws.Range("N251").Value = 32000
$ws.Range("O251").Value = 32000
$ws.Range("P251").Value = 32000
$ws.Range("Q251").Value = "$/bandeja 18 kilos granel"
$ws.Range("R251").Value = "Provincia de Limarí"
$ws.Range("S251").Value = 1778
$ws.Range("T251").Value = 18
